# Workbook: model_input_variables_jamaica_se_calibrated.xlsx
# Sheet "strategy_id-0" (sheet1.xml) is the ActiveSheet.
#
# Change: insert a new data row at row 4 for the variable
# "climate_change_factor_gnrl_hydropower_availability" (General subsector),
# which pushes the previously existing rows 4-11 down to rows 5-12.
# The new row's values: H4 = 1, I4 = 0.5, J4:AS4 = 1 (all ones).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 4; this shifts the existing
# rows 4-11 (and all their cell content/values) down to rows 5-12,
# exactly matching the rest of the diff (only row numbers change there).
$ws.Rows.Item(4).Insert()

# Populate the newly inserted row 4 with the new variable's data.
$ws.Range("A4").Value = "General"
$ws.Range("B4").Value = "climate_change_factor_gnrl_hydropower_availability"

$ws.Range("H4").Value = 1
$ws.Range("I4").Value = 0.5

# Columns J (10) through AS (45) are all 1 for this new row.
for ($col = 10; $col -le 45; $col++) {
    $ws.Cells.Item(4, $col).Value = 1
}
